$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 5: the "torch:" label in K5 was actually meant to be "N" (header for
# the existing N/C_out/H_out/W_out "torch output shape" row).
$ws.Range("K5").Value = "N"

# New rows 10-11: add the matching "input shape" header/value pair
# (N, H_in, W_in, C_in) right below the output-shape table at rows 8-9.
$ws.Range("K10").Value = "N"
$ws.Range("L10").Value = "H_in"
$ws.Range("M10").Value = "W_in"
$ws.Range("N10").Value = "C_in"

$ws.Range("K11").Value = 60
$ws.Range("L11").Value = 161
$ws.Range("M11").Value = 16
$ws.Range("N11").Value = 8

# Move the active selection to K6 (was N6).
$ws.Range("K6").Select()
